$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the score values in column C for the listed rows, keeping formatting.
$rowsToClear = @(8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 31, 32)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 3).ClearContents()
}

# Update the window view / selection to match the author's saved state.
$window = $excel.ActiveWindow
$window.ScrollRow = 9
$window.ScrollColumn = 1
$ws.Range("G12").Select()
